$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "h1-pagebreak" paragraph style from the (only, empty) first
#    paragraph of the document body so it reverts to the document default
#    ("Normal"/"Standard") paragraph style. Word omits <w:pPr><w:pStyle/></w:pPr>
#    entirely once a paragraph carries no non-default paragraph formatting.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Style = "Standard"

# ---------------------------------------------------------------------------
# 2) Add the built-in "Table Grid" table style (local/style id "Tabellenraster"
#    in this German-authored template), based on the document's existing
#    "Normal Table" ("NormaleTabelle") style, with no spacing after paragraphs
#    inside the table.
# ---------------------------------------------------------------------------
$normalTable = $d.Styles("NormaleTabelle")

$tableGrid = $d.Styles.Add("Tabellenraster", 3)
$tableGrid.NameLocal = "Table Grid"
$tableGrid.BaseStyle = $normalTable
$tableGrid.ParagraphFormat.SpaceAfter = 0

# ---------------------------------------------------------------------------
# 3) Add a custom "Table" table style (APA-style table look), also based on
#    "Normal Table", high ui priority, no spacing after paragraphs.
# ---------------------------------------------------------------------------
$tableStyle = $d.Styles.Add("Table", 3)
$tableStyle.BaseStyle = $normalTable
$tableStyle.Priority = 99
$tableStyle.ParagraphFormat.SpaceAfter = 0

Write-Output "done"
